$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.594.25"
$ws.Range("D3").Value = "'1.729.85"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'245.91"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.4821"
$ws.Range("E7").Value = "  +0.64%  "
$ws.Range("D8").Value = "'0.2670"
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("D9").Value = "'0.06181"
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("D10").Value = "'1.734.03"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").Value = "'0.07105"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "'15.63"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").Value = "'0.6129"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "'4.540"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").Value = "'77.30"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "'1.0000"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "'26.590.32"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "'0.9999"
$ws.Range("D19").Value = "'0.000006954"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").Value = "'11.55"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").Value = "'1.955.71"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").Value = "'4.523"
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("D23").Value = "'8.815"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("E24").Value = "  -1.78%  "
$ws.Range("D25").Value = "'137.43"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").Value = "'1.778"
$ws.Range("E27").Value = "  -2.21%  "
$ws.Range("D28").Value = "'1.411"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("D29").Value = "'108.19"
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("D30").Value = "'3.979"
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("D31").Value = "'0.08014"
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("D32").Value = "'3.682"
$ws.Range("E32").Value = "  -2.31%  "
$ws.Range("D33").Value = "'0.04548"
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("D34").Value = "'0.9995"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "'2.617"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("D37").Value = "'0.6329"
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("D38").Value = "'2.050"
$ws.Range("E38").Value = "  +3.53%  "
$ws.Range("D39").Value = "'0.8992"
$ws.Range("E39").Value = "  -5.22%  "
$ws.Range("D40").Value = "'2.389"
$ws.Range("E40").Value = "  -2.57%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "'102.76"
$ws.Range("E42").Value = "  -9.69%  "
$ws.Range("D43").Value = "'0.01502"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").Value = "'5.437"
$ws.Range("E44").Value = "  -4.12%  "
$ws.Range("D45").Value = "'7.157"
$ws.Range("E45").Value = "  +6.37%  "
$ws.Range("D46").Value = "'0.3899"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("D48").Value = "'0.05391"
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("D49").Value = "'7.906"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").Value = "'30.67"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").Value = "'1.255"
$ws.Range("E51").Value = "  -0.17%  "
